# "aggiornamento fino a 28/06 incluso"
# Appends the new daily rows (28 May 2021 - 28 June 2021 = serials 44344..44375)
# to the bottom of the existing data table on Sheet1 (rows 270-301),
# mirroring columns A (date), B (nuovi pos.), C (somma mobile 7gg.),
# D (somma mobile 7gg. per 100mila abitanti).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstNewRow = 270
$lastExistingRow = $firstNewRow - 1

# New rows: date serial, nuovi pos., somma mobile 7gg., somma mobile 7gg. per 100mila ab.
$data = @(
    @(44344, 1, 5, 27.81176994103905),
    @(44345, 1, 5, 27.81176994103905),
    @(44346, 0, 5, 27.81176994103905),
    @(44347, 3, 8, 44.49883190566248),
    @(44348, 2, 8, 44.49883190566248),
    @(44349, 0, 8, 44.49883190566248),
    @(44350, 0, 7, 38.93647791745467),
    @(44351, 0, 6, 33.37412392924686),
    @(44352, 0, 5, 27.81176994103905),
    @(44353, 1, 6, 33.37412392924686),
    @(44354, 1, 4, 22.24941595283124),
    @(44355, 1, 3, 16.68706196462343),
    @(44356, 0, 3, 16.68706196462343),
    @(44357, 1, 4, 22.24941595283124),
    @(44358, 0, 4, 22.24941595283124),
    @(44359, 0, 4, 22.24941595283124),
    @(44360, 0, 3, 16.68706196462343),
    @(44361, 0, 2, 11.12470797641562),
    @(44362, 0, 1, 5.56235398820781),
    @(44363, 0, 1, 5.56235398820781),
    @(44364, 0, 0, 0),
    @(44365, 0, 0, 0),
    @(44366, 0, 0, 0),
    @(44367, 0, 0, 0),
    @(44368, 1, 1, 5.56235398820781),
    @(44369, 0, 1, 5.56235398820781),
    @(44370, 2, 3, 16.68706196462343),
    @(44371, 0, 3, 16.68706196462343),
    @(44372, 1, 4, 22.24941595283124),
    @(44373, 2, 6, 33.37412392924686),
    @(44374, 0, 6, 33.37412392924686),
    @(44375, 0, 5, 27.81176994103905)
)

$lastNewRow = $firstNewRow + $data.Length - 1

# Column A carries the bordered / centered date style used throughout the
# table (same as A2:A269) - propagate it onto the freshly appended date cells.
$ws.Range("A$lastExistingRow").Copy()
$ws.Range("A${firstNewRow}:A${lastNewRow}").PasteSpecial(-4122)
$excel.CutCopyMode = 0

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $firstNewRow + $i
    $values = $data[$i]

    $ws.Cells.Item($row, 1).Value = $values[0]
    $ws.Cells.Item($row, 2).Value = $values[1]
    $ws.Cells.Item($row, 3).Value = $values[2]
    $ws.Cells.Item($row, 4).Value = $values[3]
}

Write-Host "Appended rows $firstNewRow-$lastNewRow (dates through 2021-06-28)"
